$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Alcohol" measurement sheet had a redundant column M; remove it so
# the data in column N (a duplicate/derived series) shifts left into M.
$ws.Columns("M").Delete()

# After the delete, Excel leaves the selection on the column that took M's
# place - mirror that in the saved view state.
$ws.Activate()
$ws.Range("M1").Select()
